$wb = $excel.ActiveWorkbook

# Template sheet to copy cell formatting from (same column-layout family as the
# two new sheets: TestScenario/Run/pcRegFormName/pcRegFormPcName/.../remark).
$template = $wb.Worksheets.Item("Transactions_Jewel")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------------
# Create both new worksheets (empty) first, already in their final order, so
# that sheetId/r:id allocation and sheet ordering come out right:
#   MembMang_Loan   (sheetId 22, rId19, sheet19.xml)
#   MembMang_Deposit(sheetId 23, rId20, sheet20.xml)
# ---------------------------------------------------------------------------
$wsLoan = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLoan.Name = "MembMang_Loan"

$wsDep = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLoan)
$wsDep.Name = "MembMang_Deposit"

# ---------------------------------------------------------------------------
# Fill "MembMang_Deposit" FIRST (controls the order new shared strings are
# appended in), and within its second row set column E before column A.
# ---------------------------------------------------------------------------
$wsDep.Range("A1").Value = "TestScenario"
$wsDep.Range("B1").Value = "Run"
$wsDep.Range("C1").Value = "pcRegFormName"
$wsDep.Range("D1").Value = "pcRegFormPcName"
$wsDep.Range("E1").Value = "name"
$wsDep.Range("F1").Value = "remark"

$wsDep.Range("E2").Value = "i"
$wsDep.Range("A2").Value = "MemberManagement_Deposit"
$wsDep.Range("B2").Value = "Yes"
$wsDep.Range("C2").Value = "qwerty"
$wsDep.Range("D2").Value = "zxcvb"
$wsDep.Range("F2").Value = "abcd"

# ---------------------------------------------------------------------------
# Fill "MembMang_Loan" next.
# ---------------------------------------------------------------------------
$wsLoan.Range("A1").Value = "TestScenario"
$wsLoan.Range("B1").Value = "Run"
$wsLoan.Range("C1").Value = "pcRegFormName"
$wsLoan.Range("D1").Value = "pcRegFormPcName"
$wsLoan.Range("E1").Value = "name"
$wsLoan.Range("F1").Value = "remark"

$wsLoan.Range("A2").Value = "MemberManagement_Loan"
$wsLoan.Range("B2").Value = "Yes"
$wsLoan.Range("C2").Value = "qwerty"
$wsLoan.Range("D2").Value = "zxcvb"
$wsLoan.Range("E2").Value = "n"
$wsLoan.Range("F2").Value = "abcd"

# ---------------------------------------------------------------------------
# Copy cell formatting (fill/font/border/alignment) from the existing
# "Transactions_Jewel" sheet, which uses the same header/data style pattern.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsLoan, $wsDep)) {
    $template.Range("A1:F1").Copy()
    $ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

    $template.Range("A2:D2").Copy()
    $ws.Range("A2:D2").PasteSpecial(-4122)

    $template.Range("F2").Copy()
    $ws.Range("E2:F2").PasteSpecial(-4122)

    $ws.Rows.Item(1).RowHeight = 45
    $ws.Rows.Item(2).RowHeight = 60
}

$wsDep.Range("E2:F2").NumberFormat = "General"

# ---------------------------------------------------------------------------
# Selections: MembMang_Loan -> E2, MembMang_Deposit -> H11 (which also
# becomes the active sheet / selected tab, last one touched).
# ---------------------------------------------------------------------------
$wsLoan.Activate()
$wsLoan.Range("E2").Select()

$wsDep.Activate()
$wsDep.Range("H11").Select()

# ---------------------------------------------------------------------------
# The previously-last sheet ("LoanChargePosting") is no longer the selected
# tab; its selection becomes a plain A1:D2 range (no single active cell).
# ---------------------------------------------------------------------------
$wsLCP = $wb.Worksheets.Item("LoanChargePosting")
$wsLCP.Activate()
$wsLCP.Range("A1:D2").Select()

# Leave MembMang_Deposit as the final active / selected sheet+tab.
$wsDep.Activate()

Write-Output "done"
